$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values look numeric (e.g. "50.20") but the source file stores them
# as plain text (t="inlineStr"). Prefixing with a leading apostrophe forces Excel
# to keep/treat the assigned value as text instead of auto-converting it to a number.

$ws.Range("D2").Value = "'72.063.79"
$ws.Range("E2").Value = "  +3.92%  "

$ws.Range("D3").Value = "'3.628.45"
$ws.Range("E3").Value = "  +6.51%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'594.82"
$ws.Range("E5").Value = "  +1.17%  "

$ws.Range("D6").Value = "'181.62"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").Value = "'3.612.84"
$ws.Range("E7").Value = "  +6.35%  "

$ws.Range("D8").Value = "'0.613"
$ws.Range("E8").Value = "  +2.33%  "

$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("E10").Value = "  +2.82%  "

$ws.Range("D11").Value = "'0.609"
$ws.Range("E11").Value = "  +2.34%  "

$ws.Range("D12").Value = "'50.20"
$ws.Range("E12").Value = "  +3.22%  "

$ws.Range("E13").Value = "  +1.26%  "

$ws.Range("D14").Value = "'695.89"
$ws.Range("E14").Value = "  +1.40%  "

$ws.Range("D15").Value = "'4.224.60"
$ws.Range("E15").Value = "  +6.76%  "

$ws.Range("E16").Value = "  +3.73%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'72.073.73"
$ws.Range("E17").Value = "  +3.75%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'3.578.88"
$ws.Range("E18").Value = "  +4.98%  "

$ws.Range("E19").Value = "  +2.10%  "

$ws.Range("D20").Value = "'18.43"
$ws.Range("E20").Value = "  +3.74%  "

$ws.Range("D21").Value = "'11.64"
$ws.Range("E21").Value = "  +2.69%  "

$ws.Range("D22").Value = "'0.939"
$ws.Range("E22").Value = "  +2.91%  "

$ws.Range("D23").Value = "'5.83"
$ws.Range("E23").Value = "  +7.86%  "

$ws.Range("D24").Value = "'18.01"
$ws.Range("E24").Value = "  +4.44%  "

$ws.Range("D25").Value = "'104.09"
$ws.Range("E25").Value = "  +0.55%  "

$ws.Range("D26").Value = "'4.05"
$ws.Range("E26").Value = "  +2.65%  "

$ws.Range("D27").Value = "'2.87"
$ws.Range("E27").Value = "  +5.13%  "

$ws.Range("D28").Value = "'10.05"
$ws.Range("E28").Value = "  +3.50%  "

$ws.Range("D29").Value = "'35.19"
$ws.Range("E29").Value = "  +3.13%  "

$ws.Range("D30").Value = "'9.22"
$ws.Range("E30").Value = "  +4.46%  "

$ws.Range("D31").Value = "'7.35"
$ws.Range("E31").Value = "  +5.51%  "

$ws.Range("D32").Value = "'4.20"
$ws.Range("E32").Value = "  +16.20%  "

$ws.Range("D33").Value = "'584.70"
$ws.Range("E33").Value = "  +3.77%  "

$ws.Range("D34").Value = "'11.40"
$ws.Range("E34").Value = "  +2.06%  "

$ws.Range("E35").Value = "  +3.40%  "

$ws.Range("D36").Value = "'59.64"
$ws.Range("E36").Value = "  +1.73%  "

$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("D38").Value = "'3.666.17"
$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("E39").Value = "  +1.42%  "

$ws.Range("D40").Value = "'36.11"
$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("D41").Value = "'0.0₃0771"
$ws.Range("E41").Value = "  +6.48%  "

$ws.Range("D42").Value = "'3.45"
$ws.Range("E42").Value = "  +5.26%  "

$ws.Range("D43").Value = "'0.0464"
$ws.Range("E43").Value = "  +8.70%  "

$ws.Range("E44").Value = "  +3.66%  "

$ws.Range("E45").Value = "  +3.41%  "

$ws.Range("E46").Value = "  +2.50%  "

$ws.Range("D47").Value = "'2.83"
$ws.Range("E47").Value = "  +5.63%  "

$ws.Range("E48").Value = "  +2.81%  "

$ws.Range("D49").Value = "'1.45"
$ws.Range("E49").Value = "  +3.59%  "

$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  -0.25%  "

$ws.Range("D51").Value = "'132.29"
$ws.Range("E51").Value = "  -0.27%  "
